
# -----------------------------------------------------------------------
# "minor fix for Enemy04" commit:
#   1. Colour the "Als er 2 enemies type 02 zijn ..." (laserDot) paragraph
#      and the "Enemy01 gaat niet lopen ..." paragraph green (accent 6 /
#      70AD47), both on the paragraph mark's default run formatting and on
#      the run itself.
#   2. Split the empty paragraph that used to follow the "Enemy01 ..."
#      paragraph into a brand-new paragraph that reads
#      "|--> reden: geen public static gebruiken!!" (not coloured), moving
#      the _GoBack bookmark so that it still sits right after "...geen p".
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

$laserText  = "Als er 2 enemies type 02 zijn, probleem met laserDot, niet zichtbaar"
$enemyText  = "Enemy01 gaat niet lopen als er geen Obstacle in de buurt is?!"

# Locate the two anchor paragraphs by their text (robust to absolute index
# drift) rather than hard-coding paragraph numbers.
$idxLaser = -1
$idxEnemy = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($idxLaser -eq -1 -and $t -match [regex]::Escape($laserText)) {
        $idxLaser = $i
    }
    if ($idxEnemy -eq -1 -and $t -match [regex]::Escape($enemyText)) {
        $idxEnemy = $i
    }
}

if ($idxLaser -eq -1 -or $idxEnemy -eq -1) {
    throw "Could not locate the anchor paragraphs (laser=$idxLaser, enemy=$idxEnemy)"
}

# The paragraph immediately after "Enemy01 ..." is the (currently empty)
# paragraph that gets replaced by the new "|--> reden: ..." paragraph.
$idxEmptyAfterEnemy = $idxEnemy + 1

$startPos = $d.Paragraphs.Item($idxLaser).Range.Start
$endPos   = $d.Paragraphs.Item($idxEmptyAfterEnemy).Range.End

$target = $d.Range($startPos, $endPos)

# Rebuild the whole block (4 paragraphs) in one shot so the new "reden"
# paragraph / bookmark placement come out exactly right:
#   - laserDot paragraph -> green (accent6) bold 14pt
#   - blank paragraph (unchanged, still bold 14pt, no colour)
#   - Enemy01 paragraph -> green (accent6) bold 14pt, bookmark removed
#   - new paragraph "|--> reden: geen public static gebruiken!!" (bold
#     14pt, no colour), with the _GoBack bookmark reinserted right after
#     the "p" of "public"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="28"/></w:rPr><w:t>Als er 2 enemies type 02 zijn, probleem met laserDot, niet zichtbaar</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="28"/></w:rPr><w:t>Enemy01 gaat niet lopen als er geen Obstacle in de buurt is?!</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">|--&gt; reden: geen </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>p</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>ublic static gebruiken!!</w:t></w:r></w:p>
'@

$target.InsertXML($xml)
